# Update vm_pu results for the 380 kV case (bus voltage magnitudes, rows 2-25)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.015361888628558
$ws.Cells.Item(2, 4).Value = 1.041704795992847
$ws.Cells.Item(2, 5).Value = 1.017020978142064
$ws.Cells.Item(2, 6).Value = 1.044291753189057
$ws.Cells.Item(2, 9).Value = 1.035564276191011
$ws.Cells.Item(2, 10).Value = 1.020587840857214
$ws.Cells.Item(2, 11).Value = 1.044483556421065
$ws.Cells.Item(2, 12).Value = 1.019871410994297
$ws.Cells.Item(2, 13).Value = 1.047063215003973
$ws.Cells.Item(2, 14).Value = 1.022037192654292

# Row 3
$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.016474986687234
$ws.Cells.Item(3, 4).Value = 1.042363082314531
$ws.Cells.Item(3, 5).Value = 1.017969114421507
$ws.Cells.Item(3, 6).Value = 1.045218406810522
$ws.Cells.Item(3, 9).Value = 1.035701585341233
$ws.Cells.Item(3, 10).Value = 1.021335056822347
$ws.Cells.Item(3, 11).Value = 1.044952924848954
$ws.Cells.Item(3, 12).Value = 1.020624658204479
$ws.Cells.Item(3, 13).Value = 1.047800775570264
$ws.Cells.Item(3, 14).Value = 1.022785469751803

# Row 4
$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.017195233574475
$ws.Cells.Item(4, 4).Value = 1.042786665199454
$ws.Cells.Item(4, 5).Value = 1.018583001531386
$ws.Cells.Item(4, 6).Value = 1.045815961134104
$ws.Cells.Item(4, 9).Value = 1.03578756399172
$ws.Cells.Item(4, 10).Value = 1.021818058415389
$ws.Cells.Item(4, 11).Value = 1.04525358740827
$ws.Cells.Item(4, 12).Value = 1.02111182495809
$ws.Cells.Item(4, 13).Value = 1.048275353119529
$ws.Cells.Item(4, 14).Value = 1.023269157262509

# Row 5
$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.017498025608986
$ws.Cells.Item(5, 4).Value = 1.042964169568221
$ws.Cells.Item(5, 5).Value = 1.018841170411531
$ws.Cells.Item(5, 6).Value = 1.046066680028387
$ws.Cells.Item(5, 9).Value = 1.035823021322691
$ws.Cells.Item(5, 10).Value = 1.02202099329375
$ws.Cells.Item(5, 11).Value = 1.045379253997942
$ws.Cells.Item(5, 12).Value = 1.021316573666522
$ws.Cells.Item(5, 13).Value = 1.048474223851199
$ws.Cells.Item(5, 14).Value = 1.023472380331674

# Row 6
$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.017548865758957
$ws.Cells.Item(6, 4).Value = 1.042993939858121
$ws.Cells.Item(6, 5).Value = 1.018884523419759
$ws.Cells.Item(6, 6).Value = 1.046108747881133
$ws.Cells.Item(6, 9).Value = 1.035828934381359
$ws.Cells.Item(6, 10).Value = 1.022055059975336
$ws.Cells.Item(6, 11).Value = 1.045400310998255
$ws.Cells.Item(6, 12).Value = 1.021350948578964
$ws.Cells.Item(6, 13).Value = 1.048507577445663
$ws.Cells.Item(6, 14).Value = 1.023506495391855

# Row 7
$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.017199279494635
$ws.Cells.Item(7, 4).Value = 1.0427890392633
$ws.Cells.Item(7, 5).Value = 1.018586450840242
$ws.Cells.Item(7, 6).Value = 1.045819313191293
$ws.Cells.Item(7, 9).Value = 1.035788040479538
$ws.Cells.Item(7, 10).Value = 1.021820770508741
$ws.Cells.Item(7, 11).Value = 1.045255269449767
$ws.Cells.Item(7, 12).Value = 1.021114561041246
$ws.Cells.Item(7, 13).Value = 1.048278012963522
$ws.Cells.Item(7, 14).Value = 1.023271873207344

# Row 8
$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.01573806617565
$ws.Cells.Item(8, 4).Value = 1.041927757013698
$ws.Cells.Item(8, 5).Value = 1.017341326580438
$ws.Cells.Item(8, 6).Value = 1.044605344184342
$ws.Cells.Item(8, 9).Value = 1.035611273957948
$ws.Cells.Item(8, 10).Value = 1.020840468932048
$ws.Cells.Item(8, 11).Value = 1.04464281221901
$ws.Cells.Item(8, 12).Value = 1.020126022961283
$ws.Cells.Item(8, 13).Value = 1.047313029636325
$ws.Cells.Item(8, 14).Value = 1.02229017948997

# Row 9
$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.013163173367252
$ws.Cells.Item(9, 4).Value = 1.0403919803563
$ws.Cells.Item(9, 5).Value = 1.015150169801688
$ws.Cells.Item(9, 6).Value = 1.042450522007422
$ws.Cells.Item(9, 9).Value = 1.035277859211897
$ws.Cells.Item(9, 10).Value = 1.019109243379077
$ws.Cells.Item(9, 11).Value = 1.04354029717974
$ws.Cells.Item(9, 12).Value = 1.018382296988801
$ws.Cells.Item(9, 13).Value = 1.045592199946406
$ws.Cells.Item(9, 14).Value = 1.020556495398139

# Row 10
$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.01144649612085
$ws.Cells.Item(10, 4).Value = 1.039356074893612
$ws.Cells.Item(10, 5).Value = 1.013691358524289
$ws.Cells.Item(10, 6).Value = 1.041003533014096
$ws.Cells.Item(10, 9).Value = 1.035040900391772
$ws.Cells.Item(10, 10).Value = 1.017952526347473
$ws.Cells.Item(10, 11).Value = 1.042789725306543
$ws.Cells.Item(10, 12).Value = 1.017218603765839
$ws.Cells.Item(10, 13).Value = 1.044431349294173
$ws.Cells.Item(10, 14).Value = 1.019398135695673

# Row 11
$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.010703124368803
$ws.Cells.Item(11, 4).Value = 1.038904678375257
$ws.Cells.Item(11, 5).Value = 1.013060141440473
$ws.Cells.Item(11, 6).Value = 1.040374510338154
$ws.Cells.Item(11, 9).Value = 1.034934824425467
$ws.Cells.Item(11, 10).Value = 1.017451044015562
$ws.Cells.Item(11, 11).Value = 1.042461047100394
$ws.Cells.Item(11, 12).Value = 1.016714423503887
$ws.Cells.Item(11, 13).Value = 1.043925472254111
$ws.Cells.Item(11, 14).Value = 1.01889594120133

# Row 12
$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.010426995955761
$ws.Cells.Item(12, 4).Value = 1.038736584103404
$ws.Cells.Item(12, 5).Value = 1.012825747931026
$ws.Cells.Item(12, 6).Value = 1.040140493796799
$ws.Cells.Item(12, 9).Value = 1.034894902470563
$ws.Cells.Item(12, 10).Value = 1.017264678257183
$ws.Cells.Item(12, 11).Value = 1.042338410531007
$ws.Cells.Item(12, 12).Value = 1.016527104193967
$ws.Cells.Item(12, 13).Value = 1.043737084287566
$ws.Cells.Item(12, 14).Value = 1.018709310782197

# Row 13
$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.010486226775244
$ws.Cells.Item(13, 4).Value = 1.038772660125743
$ws.Cells.Item(13, 5).Value = 1.012876023031964
$ws.Cells.Item(13, 6).Value = 1.040190707843146
$ws.Cells.Item(13, 9).Value = 1.034903489410374
$ws.Cells.Item(13, 10).Value = 1.01730465855364
$ws.Cells.Item(13, 11).Value = 1.04236474139694
$ws.Cells.Item(13, 12).Value = 1.016567286830804
$ws.Cells.Item(13, 13).Value = 1.043777515970451
$ws.Cells.Item(13, 14).Value = 1.018749347855262

# Row 14
$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.010680299658511
$ws.Cells.Item(14, 4).Value = 1.038890792321634
$ws.Cells.Item(14, 5).Value = 1.013040765001003
$ws.Cells.Item(14, 6).Value = 1.040355173989527
$ws.Cells.Item(14, 9).Value = 1.034931535079593
$ws.Cells.Item(14, 10).Value = 1.01743564086064
$ws.Cells.Item(14, 11).Value = 1.042450921152694
$ws.Cells.Item(14, 12).Value = 1.016698940536606
$ws.Cells.Item(14, 13).Value = 1.043909909886294
$ws.Cells.Item(14, 14).Value = 1.018880516172161

# Row 15
$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.010799873398642
$ws.Cells.Item(15, 4).Value = 1.03896352112863
$ws.Cells.Item(15, 5).Value = 1.013142277053138
$ws.Cells.Item(15, 6).Value = 1.040456458078847
$ws.Cells.Item(15, 9).Value = 1.034948745981392
$ws.Cells.Item(15, 10).Value = 1.017516331064394
$ws.Cells.Item(15, 11).Value = 1.042503946391926
$ws.Cells.Item(15, 12).Value = 1.016780050851081
$ws.Cells.Item(15, 13).Value = 1.043991418219265
$ws.Cells.Item(15, 14).Value = 1.018961320965259

# Row 16
$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.011495830331587
$ws.Cells.Item(16, 4).Value = 1.039385972837828
$ws.Cells.Item(16, 5).Value = 1.013733259955474
$ws.Cells.Item(16, 6).Value = 1.041045227298688
$ws.Cells.Item(16, 9).Value = 1.035047867246486
$ws.Cells.Item(16, 10).Value = 1.017985795068798
$ws.Cells.Item(16, 11).Value = 1.042811461273478
$ws.Cells.Item(16, 12).Value = 1.017252058348014
$ws.Cells.Item(16, 13).Value = 1.044464854928347
$ws.Cells.Item(16, 14).Value = 1.019431451662399

# Row 17
$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.011932374424469
$ws.Cells.Item(17, 4).Value = 1.039650205704154
$ws.Cells.Item(17, 5).Value = 1.014104090378024
$ws.Cells.Item(17, 6).Value = 1.041413887016669
$ws.Cells.Item(17, 9).Value = 1.035109114690814
$ws.Cells.Item(17, 10).Value = 1.018280112238991
$ws.Cells.Item(17, 11).Value = 1.043003374146862
$ws.Cells.Item(17, 12).Value = 1.017548057465723
$ws.Cells.Item(17, 13).Value = 1.044760967775536
$ws.Cells.Item(17, 14).Value = 1.019726186796731

# Row 18
$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.012186999506456
$ws.Cells.Item(18, 4).Value = 1.039804053890893
$ws.Cells.Item(18, 5).Value = 1.014320433674561
$ws.Cells.Item(18, 6).Value = 1.041628681901487
$ws.Cells.Item(18, 9).Value = 1.035144504306636
$ws.Cells.Item(18, 10).Value = 1.018451723029651
$ws.Cells.Item(18, 11).Value = 1.043114958944169
$ws.Cells.Item(18, 12).Value = 1.017720680448554
$ws.Cells.Item(18, 13).Value = 1.044933374525994
$ws.Cells.Item(18, 14).Value = 1.019898041294399

# Row 19
$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.012273819486611
$ws.Cells.Item(19, 4).Value = 1.039856465559965
$ws.Cells.Item(19, 5).Value = 1.014394208647224
$ws.Cells.Item(19, 6).Value = 1.041701880987541
$ws.Cells.Item(19, 9).Value = 1.035156514404632
$ws.Cells.Item(19, 10).Value = 1.018510227802512
$ws.Cells.Item(19, 11).Value = 1.043152946264847
$ws.Cells.Item(19, 12).Value = 1.01777953563231
$ws.Cells.Item(19, 13).Value = 1.044992108007026
$ws.Cells.Item(19, 14).Value = 1.019956629150747

# Row 20
$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.011885537765719
$ws.Cells.Item(20, 4).Value = 1.039621884386461
$ws.Cells.Item(20, 5).Value = 1.014064299213926
$ws.Cells.Item(20, 6).Value = 1.041374357936451
$ws.Cells.Item(20, 9).Value = 1.035102578061837
$ws.Cells.Item(20, 10).Value = 1.018248540941755
$ws.Cells.Item(20, 11).Value = 1.042982820393501
$ws.Cells.Item(20, 12).Value = 1.01751630248984
$ws.Cells.Item(20, 13).Value = 1.044729229814383
$ws.Cells.Item(20, 14).Value = 1.019694570664632

# Row 21
$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.010623150235381
$ws.Cells.Item(21, 4).Value = 1.038856017068035
$ws.Cells.Item(21, 5).Value = 1.012992250704301
$ws.Cells.Item(21, 6).Value = 1.040306753000278
$ws.Cells.Item(21, 9).Value = 1.034923290693576
$ws.Cells.Item(21, 10).Value = 1.017397072402693
$ws.Cells.Item(21, 11).Value = 1.042425558557925
$ws.Cells.Item(21, 12).Value = 1.016660173026744
$ws.Cells.Item(21, 13).Value = 1.043870936493513
$ws.Cells.Item(21, 14).Value = 1.018841892942581

# Row 22
$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.009829395746507
$ws.Cells.Item(22, 4).Value = 1.038372024779683
$ws.Cells.Item(22, 5).Value = 1.01231860818597
$ws.Cells.Item(22, 6).Value = 1.039633370036965
$ws.Cells.Item(22, 9).Value = 1.034807553897528
$ws.Cells.Item(22, 10).Value = 1.016861182732721
$ws.Cells.Item(22, 11).Value = 1.04207199952811
$ws.Cells.Item(22, 12).Value = 1.016121633957894
$ws.Cells.Item(22, 13).Value = 1.04332850117722
$ws.Cells.Item(22, 14).Value = 1.01830524224781

# Row 23
$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.01025018432808
$ws.Cells.Item(23, 4).Value = 1.038628831023549
$ws.Cells.Item(23, 5).Value = 1.012675681263341
$ws.Cells.Item(23, 6).Value = 1.039990545325519
$ws.Cells.Item(23, 9).Value = 1.034869193346684
$ws.Cells.Item(23, 10).Value = 1.017145319051569
$ws.Cells.Item(23, 11).Value = 1.042259729446253
$ws.Cells.Item(23, 12).Value = 1.016407148099718
$ws.Cells.Item(23, 13).Value = 1.043616320679596
$ws.Cells.Item(23, 14).Value = 1.018589782072821

# Row 24
$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.01190670124677
$ws.Cells.Item(24, 4).Value = 1.039634682418583
$ws.Cells.Item(24, 5).Value = 1.014082278991466
$ws.Cells.Item(24, 6).Value = 1.041392220161782
$ws.Cells.Item(24, 9).Value = 1.035105532718125
$ws.Cells.Item(24, 10).Value = 1.018262806835422
$ws.Cells.Item(24, 11).Value = 1.042992108845935
$ws.Cells.Item(24, 12).Value = 1.017530651283789
$ws.Cells.Item(24, 13).Value = 1.044743571793256
$ws.Cells.Item(24, 14).Value = 1.019708856817505

# Row 25
$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.013828855669083
$ws.Cells.Item(25, 4).Value = 1.040791147712637
$ws.Cells.Item(25, 5).Value = 1.015716291050341
$ws.Cells.Item(25, 6).Value = 1.043009441419247
$ws.Cells.Item(25, 9).Value = 1.035366647194016
$ws.Cells.Item(25, 10).Value = 1.019557258608816
$ws.Cells.Item(25, 11).Value = 1.043828073412688
$ws.Cells.Item(25, 12).Value = 1.018833304901897
$ws.Cells.Item(25, 13).Value = 1.046039484524055
$ws.Cells.Item(25, 14).Value = 1.021005146860893
